$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3338577.5
$ws.Range("J17").Value = 3338577.5
$ws.Range("L17").Value = 10015732.5
$ws.Range("N17").Value = -10016068.5

$ws.Range("H40").Value = 4972.5
$ws.Range("I40").Value = 4980
$ws.Range("K40").Value = 4980
$ws.Range("M40").Value = -4805

$ws.Range("H53").Value = 43478796
$ws.Range("I53").Value = 533.3077
$ws.Range("K53").Value = 533.3077
$ws.Range("M53").Value = 103.6923

$ws.Range("H55").Value = 269
$ws.Range("J55").Value = 75
$ws.Range("L55").Value = 75
$ws.Range("N55").Value = -503

$ws.Range("H94").Value = 1936.5
$ws.Range("I94").Value = 1658
$ws.Range("K94").Value = 1658
$ws.Range("M94").Value = -1207

$ws.Range("H132").Value = 2088.5737
$ws.Range("I132").Value = 2048.389
$ws.Range("K132").Value = 6145.167
$ws.Range("M132").Value = -3615.167

$ws.Range("H137").Value = 1962335.1
$ws.Range("I137").Value = 1045.6428
$ws.Range("J137").Value = 4349992
$ws.Range("K137").Value = 3136.9284
$ws.Range("L137").Value = 13049976
$ws.Range("M137").Value = -586.9284000000002
$ws.Range("N137").Value = -13055076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 845.9697
$ws.Range("I2").Value = 671.8889
$ws.Range("K2").Value = 671.8889
$ws.Range("M2").Value = -558.8889

$ws.Range("H116").Value = 845.9697
$ws.Range("I116").Value = 671.8889
$ws.Range("K116").Value = 671.8889
$ws.Range("M116").Value = 1622.1111

$ws.Range("H122").Value = 2801
$ws.Range("J122").Value = 4975.25
$ws.Range("L122").Value = 14925.75
$ws.Range("N122").Value = -19825.75

$ws.Range("H132").Value = 6787.294
$ws.Range("I132").Value = 6538.9
$ws.Range("K132").Value = 19616.7
$ws.Range("M132").Value = -17086.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 845.9697
$ws.Range("I3").Value = 671.8889
$ws.Range("K3").Value = 671.8889
$ws.Range("M3").Value = -557.8889

$ws.Range("H20").Value = 1590.3334
$ws.Range("I20").Value = 1353.909
$ws.Range("K20").Value = 1353.909
$ws.Range("M20").Value = -1106.909

$ws.Range("H86").Value = 4055.4285
$ws.Range("I86").Value = 2046.3334
$ws.Range("K86").Value = 2046.3334
$ws.Range("M86").Value = -923.3334

$ws.Range("H89").Value = 4055.4285
$ws.Range("I89").Value = 2046.3334
$ws.Range("K89").Value = 10231.667
$ws.Range("M89").Value = -4615.666999999999

$ws.Range("H107").Value = 8509
$ws.Range("I107").Value = 11762.16
$ws.Range("J107").Value = 3724.9412
$ws.Range("K107").Value = 11762.16
$ws.Range("L107").Value = 3724.9412
$ws.Range("M107").Value = -9842.16
$ws.Range("N107").Value = -7564.9412

$ws.Range("H140").Value = 94340.91
$ws.Range("J140").Value = 94340.91
$ws.Range("L140").Value = 94340.91
$ws.Range("N140").Value = -104700.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1412938
$ws.Range("I3").Value = 2341563.2
$ws.Range("K3").Value = 2341563.2
$ws.Range("M3").Value = -2341450.2

$ws.Range("H4").Value = 40000000
$ws.Range("J4").Value = 40000000
$ws.Range("L4").Value = 40000000
$ws.Range("N4").Value = -40000224

$ws.Range("H25").Value = 3421.1667
$ws.Range("I25").Value = 3421.1667
$ws.Range("K25").Value = 3421.1667
$ws.Range("M25").Value = -3247.1667

$ws.Range("H31").Value = 4113.8887
$ws.Range("I31").Value = 1729.5
$ws.Range("J31").Value = 4980.9395
$ws.Range("K31").Value = 1729.5
$ws.Range("L31").Value = 4980.9395
$ws.Range("M31").Value = -1434.5
$ws.Range("N31").Value = -5570.9395

$ws.Range("H34").Value = 4113.8887
$ws.Range("I34").Value = 1729.5
$ws.Range("J34").Value = 4980.9395
$ws.Range("K34").Value = 1729.5
$ws.Range("L34").Value = 4980.9395
$ws.Range("M34").Value = -1527.5
$ws.Range("N34").Value = -5384.9395

$ws.Range("H59").Value = 52500
$ws.Range("J59").Value = 52500
$ws.Range("L59").Value = 52500
$ws.Range("N59").Value = -54790

$ws.Range("H60").Value = 30999.8
$ws.Range("I60").Value = 18666.334
$ws.Range("J60").Value = 49500
$ws.Range("K60").Value = 18666.334
$ws.Range("L60").Value = 49500
$ws.Range("M60").Value = -18155.334
$ws.Range("N60").Value = -50522

$ws.Range("H62").Value = 351666.34
$ws.Range("J62").Value = 502499.5
$ws.Range("L62").Value = 502499.5
$ws.Range("N62").Value = -503747.5

$ws.Range("H65").Value = 351666.34
$ws.Range("J65").Value = 502499.5
$ws.Range("L65").Value = 2512497.5
$ws.Range("N65").Value = -2518737.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1021.8947
$ws.Range("I113").Value = 815.8
$ws.Range("J113").Value = 1095.5
$ws.Range("K113").Value = 2447.4
$ws.Range("L113").Value = 3286.5
$ws.Range("M113").Value = -277.3999999999996
$ws.Range("N113").Value = -7626.5

$ws.Range("H121").Value = 2694.75
$ws.Range("J121").Value = 2694.75
$ws.Range("L121").Value = 8084.25
$ws.Range("N121").Value = -10704.25

$ws.Range("H131").Value = 7812557
$ws.Range("J131").Value = 8629143
$ws.Range("L131").Value = 25887429
$ws.Range("N131").Value = -25897509

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6978.103
$ws.Range("I70").Value = 7398.9443
$ws.Range("J70").Value = 5354.857
$ws.Range("K70").Value = 7398.9443
$ws.Range("L70").Value = 5354.857
$ws.Range("M70").Value = -7128.9443
$ws.Range("N70").Value = -5894.857

$ws.Range("H73").Value = 6978.103
$ws.Range("I73").Value = 7398.9443
$ws.Range("J73").Value = 5354.857
$ws.Range("K73").Value = 7398.9443
$ws.Range("L73").Value = 5354.857
$ws.Range("M73").Value = -6462.9443
$ws.Range("N73").Value = -7226.857

$ws.Range("H80").Value = 14337215
$ws.Range("I80").Value = 58473.76
$ws.Range("K80").Value = 58473.76
$ws.Range("M80").Value = -57475.76

$ws.Range("H83").Value = 14337215
$ws.Range("I83").Value = 58473.76
$ws.Range("K83").Value = 292368.8
$ws.Range("M83").Value = -287376.8

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 8617733
$ws.Range("I132").Value = 4002.8948
$ws.Range("J132").Value = 21207032
$ws.Range("K132").Value = 12008.6844
$ws.Range("L132").Value = 63621096
$ws.Range("M132").Value = -9478.6844
$ws.Range("N132").Value = -63626156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3218.7827
$ws.Range("I40").Value = 2953.9048
$ws.Range("K40").Value = 2953.9048
$ws.Range("M40").Value = -2817.9048

$ws.Range("H122").Value = 4469.5483
$ws.Range("I122").Value = 3588.7
$ws.Range("J122").Value = 6071.091
$ws.Range("K122").Value = 10766.1
$ws.Range("L122").Value = 18213.273
$ws.Range("M122").Value = -8316.099999999999
$ws.Range("N122").Value = -23113.273

$ws.Range("H132").Value = 11130.923
$ws.Range("I132").Value = 3333.3333
$ws.Range("J132").Value = 17814.572
$ws.Range("K132").Value = 9999.999899999999
$ws.Range("L132").Value = 53443.716
$ws.Range("M132").Value = -7469.999899999999
$ws.Range("N132").Value = -58503.716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 3500
$ws.Range("I52").Value = 3500
$ws.Range("K52").Value = 3500
$ws.Range("M52").Value = -3274

$ws.Range("H61").Value = 10182.833
$ws.Range("I61").Value = 8219.4
$ws.Range("K61").Value = 8219.4
$ws.Range("M61").Value = -7927.4

$ws.Range("H94").Value = 19996.5
$ws.Range("J94").Value = 19996.5
$ws.Range("L94").Value = 19996.5
$ws.Range("N94").Value = -21798.5

$ws.Range("H113").Value = 720.8214
$ws.Range("I113").Value = 838.8946999999999
$ws.Range("K113").Value = 2516.6841
$ws.Range("M113").Value = -346.6840999999999

$ws.Range("H136").Value = 2968.2354
$ws.Range("I136").Value = 1565.8422
$ws.Range("J136").Value = 4744.6
$ws.Range("K136").Value = 4697.5266
$ws.Range("L136").Value = 14233.8
$ws.Range("M136").Value = -2147.5266
$ws.Range("N136").Value = -19333.8
